$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "758×9=" "212×8="
Replace-Text "759×9=" "424×9="
Replace-Text "222×5=" "616×3="
Replace-Text "609×5=" "926×3="
Replace-Text "476×6=" "207×4="
Replace-Text "360×6=" "582×9="
Replace-Text "832×6=" "935×8="
Replace-Text "776×9=" "817×5="
Replace-Text "177×6=" "341×2="
Replace-Text "600×6=" "967×4="
Replace-Text "526×2=" "153×2="
Replace-Text "367×9=" "318×3="
Replace-Text "903×6=" "390×4="
Replace-Text "562×8=" "418×8="
Replace-Text "236×3=" "422×8="
Replace-Text "264×8=" "429×7="
Replace-Text "681×9=" "186×2="
Replace-Text "149×8=" "786×4="
Replace-Text "497×2=" "189×6="
Replace-Text "840×8=" "333×6="
Replace-Text "199×7=" "726×8="
Replace-Text "359×7=" "940×5="
Replace-Text "305×2=" "871×5="
Replace-Text "863×9=" "919×4="
Replace-Text "486×2=" "457×2="
